$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new "2009" row of IPCA monthly data as row 17.
$ws.Range("A17").Value = "`"2009`""
$ws.Range("B17").Value = "`"0,48%"
$ws.Range("C17").Value = "`"0,55%`""
$ws.Range("D17").Value = "`"0,20%`""
$ws.Range("E17").Value = "`"0,48%`""
$ws.Range("F17").Value = "`"0,47%`""
$ws.Range("G17").Value = "`"0,36%`""
$ws.Range("H17").Value = "`"0,24%`""
$ws.Range("I17").Value = "`"0,15%`""
$ws.Range("J17").Value = "`"0,24%`""
$ws.Range("K17").Value = "`"0,28%`""
$ws.Range("L17").Value = "`"0,41%`""
$ws.Range("M17").Value = "`"0,37%`""
$ws.Range("N17").Value = "`"4,31%`""

# Update the selection to reflect the cell the author ended up on.
$ws.Range("N18").Select()
